$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Ankir"
$ws.Range("B5").Value = "Male"
